$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "RUNMANAGER" (sheet1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("RUNMANAGER")

# Delete the old row 6 (enterPIMPageTest) - its content is being consolidated
# into rows 3-5 below, so the table shrinks from 6 to 5 data-ish rows.
$ws1.Rows.Item(6).Delete()

# Row 3: newTest -> enterAdminPageTest
$ws1.Cells.Item(3, 1).Value = "enterAdminPageTest"
$ws1.Cells.Item(3, 2).Value = "To check whether use can enter Admin page"
$ws1.Cells.Item(3, 4).Value = "1"

# Row 4: googleSearchTest -> enterPIMPageTest (also pick up the text style
# used by the rest of the table, since this row previously had column A
# left with the default/no style)
$fmtSrc = $ws1.Cells.Item(3, 1)
$fmtSrc.Copy()
$ws1.Cells.Item(4, 1).PasteSpecial(-4122)
$ws1.Cells.Item(4, 1).Value = "enterPIMPageTest"
$ws1.Cells.Item(4, 2).Value = "To check whether use can enter PIM page"
$ws1.Cells.Item(4, 4).Value = "1"

# Row 5: enterAdminPageTest -> baiduSearchTest
$ws1.Cells.Item(5, 1).Value = "baiduSearchTest"
$ws1.Cells.Item(5, 2).Value = "To check whether baidu search is working"

$ws1.Range("E11").Select()

# ---------------------------------------------------------------------------
# Sheet "DATA" (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DATA")

# Insert a new "version" column before the "username" column (D).
$ws2.Columns.Item(4).Insert()

# Give the new column D the same "quote-prefixed text" style already used
# by the old last column (now shifted to H) so version numbers like
# "116.0" are stored/displayed as text.
$verFmtSrc = $ws2.Cells.Item(2, 8)
$verFmtSrc.Copy()
$ws2.Range("D2:D10").PasteSpecial(-4122)

# Header row
$ws2.Cells.Item(1, 4).Value = "version"

# Row 2
$ws2.Cells.Item(2, 2).Value = "no"
$ws2.Cells.Item(2, 4).Value = "116.0"
$ws2.Cells.Item(2, 6).Value = "YWRtaW4xMjM="

# Row 3
$ws2.Cells.Item(3, 2).Value = "no"
$ws2.Cells.Item(3, 3).Value = "firefox"
$ws2.Cells.Item(3, 4).Value = "116.0"
$ws2.Cells.Item(3, 6).Value = "YWRtaW4xMjM="

# Row 4
$ws2.Cells.Item(4, 4).Value = ""
$ws2.Cells.Item(4, 6).Value = "YWRtaW4xMjM="

# Row 5
$ws2.Cells.Item(5, 4).Value = "116.0"
$ws2.Cells.Item(5, 6).Value = "YWRtaW4xMjM="

# Row 6
$ws2.Cells.Item(6, 2).Value = "no"
$ws2.Cells.Item(6, 3).Value = "chrome"
$ws2.Cells.Item(6, 4).Value = "116.0"
$ws2.Cells.Item(6, 6).Value = "YWRtaW4xMjM="

# Row 7
$ws2.Cells.Item(7, 1).Value = "baiduSearchTest"
$ws2.Cells.Item(7, 2).Value = "yes"
$ws2.Cells.Item(7, 3).Value = "chrome"
$ws2.Cells.Item(7, 4).Value = "115.0"

# Row 8
$ws2.Cells.Item(8, 1).Value = "baiduSearchTest"
$ws2.Cells.Item(8, 2).Value = "yes"
$ws2.Cells.Item(8, 3).Value = "firefox"
$ws2.Cells.Item(8, 4).Value = "116.0"
$ws2.Cells.Item(8, 8).Value = "Hello world"

# Row 9
$ws2.Cells.Item(9, 2).Value = "no"
$ws2.Cells.Item(9, 4).Value = "116.0"
$ws2.Cells.Item(9, 6).Value = "YWRtaW4xMjM="

# Row 10
$ws2.Cells.Item(10, 2).Value = "no"
$ws2.Cells.Item(10, 3).Value = "firefox"
$ws2.Cells.Item(10, 4).Value = "116.0"
$ws2.Cells.Item(10, 6).Value = "YWRtaW4xMjM="

# Column width for the (new) password column F, which now holds the longer
# base64-encoded values.
$ws2.Columns.Item(6).ColumnWidth = 17.1796875

$ws2.Range("D14").Select()
